$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 39, last existing data row) into row 40
# for the columns that should keep the "Arial, centered" look (style used across the table).
$ws.Range("A39:C39").Copy()
$ws.Range("A40:C40").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("E39:G39").Copy()
$ws.Range("E40:G40").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("I39:L39").Copy()
$ws.Range("I40:L40").PasteSpecial(-4122) # xlPasteFormats

# D40 and H40 keep the default font but centered alignment (matches target style).
$ws.Range("D40").HorizontalAlignment = -4108 # xlCenter
$ws.Range("H40").HorizontalAlignment = -4108 # xlCenter

# New data row describing the "FTV 30 20 SGW" product.
$ws.Range("A40").Value = "FTV 30 20 SGW"
$ws.Range("B40").Value = "Trina Solar S+ 470 W"
$ws.Range("C40").Value = "Monofase"
$ws.Range("D40").Value = "30 Kw"
$ws.Range("E40").Value = "25 Anni"
$ws.Range("F40").Value = "Sungrow"
$ws.Range("G40").Value = "Sungrow"
$ws.Range("H40").Value = "20 Kw"
$ws.Range("I40").Value = "10 Anni"
$ws.Range("J40").Formula = "=L40+K40"
$ws.Range("K40").Formula = "=L40*0.06"
$ws.Range("L40").Value = 39336

# Row 41 gains an (empty, but styled) J cell, matching the shared-formula fill-down.
$ws.Range("J40").Copy()
$ws.Range("J41").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# Restore the view state recorded for this sheet.
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Range("J43").Select()
